# Update the "想去人数" (interested-people count) values in column F
# for both the "展览" and "全部类型" sheets, reflecting new figures
# from the regenerated GitHub Pages data (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Map of row -> new value for column F on each affected worksheet.
$updates = @{
    2  = 1904
    7  = 1580
    19 = 3668
    23 = 331
    24 = 589
    25 = 321
    28 = 1465
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
